# Replace the worksheet's 25 "two-digit divided by one-digit" answer
# cells (data rows 1, 5, 9, 13, 17 of the single 20x5 table; blank rows
# in between are left untouched) with the freshly generated problems.
$d = $word.ActiveDocument
$table = $d.Tables(1)

$updates = @(
    @{Row=1; Col=1; Text="93÷2=46, 1"}
    @{Row=1; Col=2; Text="16÷5=3, 1"}
    @{Row=1; Col=3; Text="23÷5=4, 3"}
    @{Row=1; Col=4; Text="50÷8=6, 2"}
    @{Row=1; Col=5; Text="93÷7=13, 2"}
    @{Row=5; Col=1; Text="67÷9=7, 4"}
    @{Row=5; Col=2; Text="31÷6=5, 1"}
    @{Row=5; Col=3; Text="98÷8=12, 2"}
    @{Row=5; Col=4; Text="47÷5=9, 2"}
    @{Row=5; Col=5; Text="92÷8=11, 4"}
    @{Row=9; Col=1; Text="44÷5=8, 4"}
    @{Row=9; Col=2; Text="83÷5=16, 3"}
    @{Row=9; Col=3; Text="49÷2=24, 1"}
    @{Row=9; Col=4; Text="48÷8=6, 0"}
    @{Row=9; Col=5; Text="28÷2=14, 0"}
    @{Row=13; Col=1; Text="25÷2=12, 1"}
    @{Row=13; Col=2; Text="62÷3=20, 2"}
    @{Row=13; Col=3; Text="30÷6=5, 0"}
    @{Row=13; Col=4; Text="75÷4=18, 3"}
    @{Row=13; Col=5; Text="94÷6=15, 4"}
    @{Row=17; Col=1; Text="91÷3=30, 1"}
    @{Row=17; Col=2; Text="40÷6=6, 4"}
    @{Row=17; Col=3; Text="33÷8=4, 1"}
    @{Row=17; Col=4; Text="58÷3=19, 1"}
    @{Row=17; Col=5; Text="89÷4=22, 1"}
)

foreach ($u in $updates) {
    $cell = $table.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
